# Line_Data.xlsx update:
# Insert two new leading columns (node1/node2 pair ids) and shift the
# existing four metric columns (eff. impedance, reactance, resistance,
# rating) one step to the right so the sheet becomes A:F instead of A:D.
# Also re-applies the default "General" number format to the new A1 cell,
# switches the page to portrait orientation, and leaves the selection on I9
# - matching the state the workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# node1, node2, impedance, reactance, resistance, (rating - only rows 1-2)
$rows = @(
  @(1, 2, 0.019380000000000001, 0.05917, 0.0528, 95),
  @(1, 5, 0.5403, 0.22303999999999999, 0.049200000000000001, 100),
  @(2, 3, 0.046989999999999997, 0.19797000000000001, 0.043799999999999999),
  @(2, 4, 0.058110000000000002, 0.17632, 0.034000000000000002),
  @(2, 5, 0.56950000000000001, 0.17388000000000001, 0.034599999999999999),
  @(3, 4, 0.67010000000000003, 0.17102999999999999, 0.012800000000000001),
  @(4, 5, 0.013350000000000001, 0.042110000000000002, 0),
  @(4, 7, 0, 0.55618000000000001, 0),
  @(5, 6, 0, 0.25202000000000002, 0),
  @(6, 9, 0.094979999999999995, 0.19889999999999999, 0),
  @(6, 10, 0.12291000000000001, 0.25580999999999998, 0),
  @(6, 11, 0.06615, 0.13027, 0),
  @(7, 8, 0.031809999999999998, 0.084500000000000006, 0),
  @(7, 12, 0.12711, 0.27038000000000001, 0),
  @(8, 9, 0.082049999999999998, 0.19206999999999999, 0),
  @(10, 11, 0.22092000000000001, 0.19988, 0),
  @(11, 12, 0.17093, 0.34802, 0)
)

$r = 1
foreach ($rowVals in $rows) {
  $c = 1
  foreach ($v in $rowVals) {
    $ws.Cells.Item($r, $c).Value = $v
    $c = $c + 1
  }
  $r = $r + 1
}

# A1 re-applies General formatting explicitly (flags applyNumberFormat on its style).
$ws.Range("A1").NumberFormat = "General"

# Page orientation flipped to portrait.
$ws.PageSetup.Orientation = 1

# Final UI selection left on I9.
$ws.Range("I9").Select()
